$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$aValues = @(
    5390,
    5340,
    5300,
    5260,
    5230,
    5210,
    5190,
    5170,
    5150,
    5150,
    5150,
    5160,
    5190,
    5200,
    5200,
    5210,
    5250,
    5290,
    5340,
    5400,
    5470,
    5560,
    5650,
    5750,
    5900,
    6000,
    6050,
    6060,
    6040,
    6010,
    5960,
    5860,
    5720,
    5610,
    5490,
    5380,
    5310,
    5220,
    5160,
    5100,
    5030,
    4990,
    4960,
    4940,
    4910,
    4890,
    4870,
    4860,
    4870,
    4880,
    4880,
    4880,
    4890,
    4900,
    4910,
    4910,
    4930,
    4940,
    4960,
    4980,
    5040,
    5120,
    5180,
    5250,
    5320,
    5410,
    5500,
    5600,
    5740,
    5860,
    5990,
    6130,
    6260,
    6400,
    6520,
    6650,
    6740,
    6830,
    6920,
    7040,
    7120,
    7140,
    7120,
    7000,
    6790,
    6630,
    6480,
    6300,
    6100,
    5940,
    5830,
    5720,
    5640,
    5600,
    5540,
    5460
)

$bValues = @(
    45797,
    45797.01041666666,
    45797.02083333334,
    45797.03125,
    45797.04166666666,
    45797.05208333334,
    45797.0625,
    45797.07291666666,
    45797.08333333334,
    45797.09375,
    45797.10416666666,
    45797.11458333334,
    45797.125,
    45797.13541666666,
    45797.14583333334,
    45797.15625,
    45797.16666666666,
    45797.17708333334,
    45797.1875,
    45797.19791666666,
    45797.20833333334,
    45797.21875,
    45797.22916666666,
    45797.23958333334,
    45797.25,
    45797.26041666666,
    45797.27083333334,
    45797.28125,
    45797.29166666666,
    45797.30208333334,
    45797.3125,
    45797.32291666666,
    45797.33333333334,
    45797.34375,
    45797.35416666666,
    45797.36458333334,
    45797.375,
    45797.38541666666,
    45797.39583333334,
    45797.40625,
    45797.41666666666,
    45797.42708333334,
    45797.4375,
    45797.44791666666,
    45797.45833333334,
    45797.46875,
    45797.47916666666,
    45797.48958333334,
    45797.5,
    45797.51041666666,
    45797.52083333334,
    45797.53125,
    45797.54166666666,
    45797.55208333334,
    45797.5625,
    45797.57291666666,
    45797.58333333334,
    45797.59375,
    45797.60416666666,
    45797.61458333334,
    45797.625,
    45797.63541666666,
    45797.64583333334,
    45797.65625,
    45797.66666666666,
    45797.67708333334,
    45797.6875,
    45797.69791666666,
    45797.70833333334,
    45797.71875,
    45797.72916666666,
    45797.73958333334,
    45797.75,
    45797.76041666666,
    45797.77083333334,
    45797.78125,
    45797.79166666666,
    45797.80208333334,
    45797.8125,
    45797.82291666666,
    45797.83333333334,
    45797.84375,
    45797.85416666666,
    45797.86458333334,
    45797.875,
    45797.88541666666,
    45797.89583333334,
    45797.90625,
    45797.91666666666,
    45797.92708333334,
    45797.9375,
    45797.94791666666,
    45797.95833333334,
    45797.96875,
    45797.97916666666,
    45797.98958333334
)

$startRow = 2
for ($i = 0; $i -lt $aValues.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $aValues[$i]
    $ws.Cells.Item($startRow + $i, 2).Value = $bValues[$i]
}
